# Apply the "Add files via upload" change to bd/PRIORIZACION_ACTIVOS_R_V.xlsx
#
# Summary of the change:
#  - Column B ("Observacion") previously held either "x" (marking an
#    inactive asset) or was left blank (active asset).
#  - The "x" marker is renamed to "inactivo", and every previously-blank
#    row is now explicitly marked "activo".
#  - The sheet's zoom was bumped from 150% to 160%, and the oversized
#    full-column selection on B was tightened to just B1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows that already show "x" in column B: rename the label to "inactivo" ---
$inactiveRows = @(15, 20, 27, 28, 31, 32, 34, 36, 37, 38, 39, 40)
foreach ($r in $inactiveRows) {
    $ws.Range("B$r").Value = "inactivo"
}

# --- Rows B2:B11 already carry the row's style (s="14") but are blank:
#     just fill in "activo" ---
$activeStyledRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11)
foreach ($r in $activeStyledRows) {
    $ws.Range("B$r").Value = "activo"
}

# --- Rows B12:B35 have no cell at all yet: pick up the same style other
#     "Observacion" cells use (copy format from B2) before filling the
#     value in ---
$activeUnstyledRows = @(12, 13, 14, 16, 17, 18, 19, 21, 22, 23, 24, 25, 26, 29, 30, 33, 35)
$ws.Range("B2").Copy()
foreach ($r in $activeUnstyledRows) {
    $ws.Range("B$r").PasteSpecial(-4122)
}
foreach ($r in $activeUnstyledRows) {
    $ws.Range("B$r").Value = "activo"
}

# --- View tweaks: zoom 150% -> 160%, and shrink the whole-column
#     selection down to B1 ---
$ws.Application.ActiveWindow.Zoom = 160
$ws.Range("B1").Select()
